$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the full B2:D9 block to 0 first
$ws.Range("B2:D9").Value = 0

# Then apply the specific non-zero overrides from the diff
$ws.Range("D5").Value = -0.6965100207141295
$ws.Range("D9").Value = -0.6324894126781301
